$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 400
$ws.Range("B2").Value = 375
$ws.Range("C2").Value = 9
$ws.Range("D2").Value = 16

$ws.Range("B5").Value = 0.9375
$ws.Range("C5").Value = 0.0225
$ws.Range("D5").Value = 0.04
